$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly observation was added to the series. It belongs right before
# the existing row 191 (chronologically it is the most recent date), so shift
# every row from 191 downward by one and populate the freshly inserted row
# with the new record.
$ws.Rows(191).Insert()

$ws.Cells.Item(191, 1).Value = 10
$ws.Cells.Item(191, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(191, 3).Value = "La Araucanía"
$ws.Cells.Item(191, 4).Value = 44719
$ws.Cells.Item(191, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(191, 5).Value = 9
$ws.Cells.Item(191, 6).Value = 100112044
$ws.Cells.Item(191, 7).Value = "Perejil"
$ws.Cells.Item(191, 8).Value = "Sin especificar"
$ws.Cells.Item(191, 9).Value = "Primera"
$ws.Cells.Item(191, 10).Value = 20
$ws.Cells.Item(191, 11).Value = 4000
$ws.Cells.Item(191, 12).Value = 4000
$ws.Cells.Item(191, 13).Value = 4000
$ws.Cells.Item(191, 14).Value = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(191, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(191, 16).Value = 1333
$ws.Cells.Item(191, 17).Value = 3
$ws.Cells.Item(191, 18).Value = "Hortaliza"
